# Add two new columns, I ("I0") and J ("IF"), to the pitching log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), matching the bold/centered/bordered style already
# used by the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# I0 / IF values for data rows 2..76, in row order.
$values = @(
    @(7,7),
    @(10,10),
    @(9,9),
    @(7,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,7),
    @(7,7),
    @(8,8),
    @(5,5),
    @(9,9),
    @(6,7),
    @(6,6),
    @(6,6),
    @(7,7),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(6,6),
    @(6,6),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(10,10),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,7),
    @(8,8),
    @(9,9),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,6),
    @(6,6),
    @(8,8),
    @(8,8),
    @(5,5)
)

for ($idx = 0; $idx -lt $values.Count; $idx++) {
    $row = $idx + 2
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
